$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update "최종점수" (K column) values - decreased by 0.5
$ws.Range("K2").Value = 54.1
$ws.Range("K3").Value = 50.1
$ws.Range("K4").Value = 45.9
$ws.Range("K5").Value = 44.7
$ws.Range("K6").Value = 36.1

# Update "MACRO_SCORE" (N column) values - new constant for all rows
$ws.Range("N2").Value = 49.16024380385575
$ws.Range("N3").Value = 49.16024380385575
$ws.Range("N4").Value = 49.16024380385575
$ws.Range("N5").Value = 49.16024380385575
$ws.Range("N6").Value = 49.16024380385575
